$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from Sheet1 to ProfOffices
$ws.Name = "ProfOffices"

# Fill Name (B) column first for both rows, matching the shared-string order
$ws.Range("B2").Value = "ΤΕΣΤ1"
$ws.Range("B3").Value = "ΤΕΣΤ2"

# Then the rest of row 2 (ID=1 already present in A2)
$ws.Range("C2").Value = "ΤΣΕΤ1"
$ws.Range("D2").Value = "asd@asd.com"
$ws.Range("E2").Value = "Β.098"
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 2154002488

# Then the rest of row 3 (ID=2 already present in A3)
$ws.Range("C3").Value = "ΤΣΕΤ2"
$ws.Range("D3").Value = "aqweqe@vdfjvio.gr"
$ws.Range("E3").Value = "Α.104"
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2541896302

# Add hyperlinks for the mail addresses (mailto:)
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:asd@asd.com", "", "", "asd@asd.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:aqweqe@vdfjvio.gr", "", "", "aqweqe@vdfjvio.gr")

# Reset the saved selection back to the default top-left cell (the recorded
# selection in the source file had drifted to H17; the updated workbook
# no longer carries that stale UI state)
$ws.Range("A1").Select() | Out-Null
